# FN-3481 cherry-pick: add "Crumpet" and "Scone" GEF/exporter rows to the
# utilisation report fixture (rows 5 and 6), fixing the numerical rounding
# error covered by #3546.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 already carries the (previously-empty) formatting for columns A:K,
# so borrow that formatting for the brand-new row 5 (columns A:J only).
$ws.Range("A6:J6").Copy()
$ws.Range("A5:J5").PasteSpecial(-4122)  # xlPasteFormats

# New row 5 - "Crumpet" facility.
$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# Existing (blank) row 6 - "Scone" facility.
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# Match the author's final selection in the saved workbook.
$ws.Range("A5:J6").Select()
